# Updated cryptos list on Tue Jan 30 05:56:21 UTC 2024 with GitHub Actions
#
# Refreshes the price / 1h-volume columns (and, for rows 36-37, the coin
# name + link as well - Celestia and LidoDAOToken swap rank) to match the
# latest coinranking.com snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    # Writes $Value into $Cell as literal TEXT, even when the string looks
    # like a number (e.g. "311.29"). Excel normally auto-converts such
    # strings to numbers on assignment, which would change the cell type
    # from the original inline-string/text cells this sheet uses - so we
    # force a text NumberFormat for the write, then restore the default
    # "Normal" style (the cells carry no special formatting otherwise).
    param($Cell, $Value)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Value
    $Cell.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = '43.363.66'
$ws.Range("E2").Value = '  +2.89%  '

# Row 3
$ws.Range("D3").Value = '2.306.30'
$ws.Range("E3").Value = '  +1.87%  '

# Row 4
$ws.Range("E4").Value = '  -0.05%  '

# Row 5
Set-TextValue $ws.Range("D5") '311.29'
$ws.Range("E5").Value = '  +1.62%  '

# Row 6
Set-TextValue $ws.Range("D6") '102.67'
$ws.Range("E6").Value = '  +6.63%  '

# Row 7
$ws.Range("E7").Value = '  +1.48%  '

# Row 8
$ws.Range("E8").Value = '  -0.03%  '

# Row 9
$ws.Range("E9").Value = '  +7.73%  '

# Row 10
Set-TextValue $ws.Range("D10") '35.83'
$ws.Range("E10").Value = '  +2.70%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.0813'
$ws.Range("E11").Value = '  +3.09%  '

# Row 12
$ws.Range("E12").Value = '  -0.79%  '

# Row 13
$ws.Range("E13").Value = '  +1.35%  '

# Row 14
$ws.Range("D14").Value = '2.663.51'

# Row 15
Set-TextValue $ws.Range("D15") '15.03'
$ws.Range("E15").Value = '  +2.62%  '

# Row 16
$ws.Range("D16").Value = '2.323.84'
$ws.Range("E16").Value = '  +2.83%  '

# Row 17
$ws.Range("E17").Value = '  +2.38%  '

# Row 18
$ws.Range("D18").Value = '43.269.91'
$ws.Range("E18").Value = '  +2.98%  '

# Row 19
$ws.Range("E19").Value = '  +0.80%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0933'
$ws.Range("E20").Value = '  +3.25%  '

# Row 21
$ws.Range("E21").Value = '  +3.24%  '

# Row 22
Set-TextValue $ws.Range("D22") '68.06'
$ws.Range("E22").Value = '  +0.49%  '

# Row 23
Set-TextValue $ws.Range("D23") '241.42'
$ws.Range("E23").Value = '  +1.94%  '

# Row 24
$ws.Range("E24").Value = '  +1.71%  '

# Row 25
$ws.Range("E25").Value = '  +2.86%  '

# Row 26
$ws.Range("E26").Value = '  +0.06%  '

# Row 27
Set-TextValue $ws.Range("D27") '24.62'
$ws.Range("E27").Value = '  +4.94%  '

# Row 28
$ws.Range("E28").Value = '  +8.54%  '

# Row 29
Set-TextValue $ws.Range("D29") '36.96'
$ws.Range("E29").Value = '  -1.37%  '

# Row 30
Set-TextValue $ws.Range("D30") '9.64'
$ws.Range("E30").Value = '  +1.03%  '

# Row 31
Set-TextValue $ws.Range("D31") '168.77'
$ws.Range("E31").Value = '  +3.89%  '

# Row 32
Set-TextValue $ws.Range("D32") '5.29'
$ws.Range("E32").Value = '  +1.17%  '

# Row 33
$ws.Range("E33").Value = '  +0.07%  '

# Row 34
Set-TextValue $ws.Range("D34") '2.52'
$ws.Range("E34").Value = '  +6.35%  '

# Row 35
$ws.Range("E35").Value = '  +1.06%  '

# Row 36
$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range("D36") '3.08'
$ws.Range("E36").Value = '  -2.98%  '

# Row 37
$ws.Range("B37").Value = 'Celestia'
$ws.Range("C37").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue $ws.Range("D37") '17.64'
$ws.Range("E37").Value = '  +0.22%  '

# Row 38
Set-TextValue $ws.Range("D38") '1.90'
$ws.Range("E38").Value = '  +4.47%  '

# Row 39
$ws.Range("E39").Value = '  +1.42%  '

# Row 40
$ws.Range("E40").Value = '  +1.76%  '

# Row 41
Set-TextValue $ws.Range("D41") '4.37'
$ws.Range("E41").Value = '  +7.77%  '

# Row 42
Set-TextValue $ws.Range("D42") '2.31'
$ws.Range("E42").Value = '  -0.84%  '

# Row 43
$ws.Range("E43").Value = '  +3.63%  '

# Row 44
Set-TextValue $ws.Range("D44") '0.0289'
$ws.Range("E44").Value = '  +3.22%  '

# Row 45
$ws.Range("D45").Value = '1.971.19'
$ws.Range("E45").Value = '  +1.20%  '

# Row 46
Set-TextValue $ws.Range("D46") '2.98'
$ws.Range("E46").Value = '  +2.49%  '

# Row 47
Set-TextValue $ws.Range("D47") '9.87'
$ws.Range("E47").Value = '  +0.37%  '

# Row 48
Set-TextValue $ws.Range("D48") '55.43'
$ws.Range("E48").Value = '  +2.88%  '

# Row 49
Set-TextValue $ws.Range("D49") '2.90'
$ws.Range("E49").Value = '  +0.46%  '

# Row 50
$ws.Range("E50").Value = '  +8.09%  '

# Row 51
$ws.Range("D51").Value = '2.532.74'
$ws.Range("E51").Value = '  +1.76%  '
